$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width (21.83203125 -> ~23.33203125) ---
$ws.Columns.Item(1).ColumnWidth = 22.5

# --- Row 24 ---
$ws.Range("A24").Value = "5 - Change citizen's trashes "

$ws.Range("B24:B25").Merge()
$ws.Range("B24").Value = "Team managers"
$ws.Range("B24").WrapText = $true
$ws.Range("B24").VerticalAlignment = -4108
$ws.Range("B24").HorizontalAlignment = -4131

$ws.Range("C24").Value = "Get information from trashes sensors"
$ws.Range("C24").WrapText = $true
$ws.Range("C24").VerticalAlignment = -4108

$ws.Range("D24").Value = "1h"
$ws.Range("D24").WrapText = $true
$ws.Range("D24").VerticalAlignment = -4108
$ws.Range("D24").HorizontalAlignment = -4108

$ws.Range("E24").Value = "Team manager "
$ws.Range("E24").WrapText = $true
$ws.Range("E24").VerticalAlignment = -4108

# --- Row 25 ---
$ws.Range("C25").Value = "Adapte garbage collector's route"
$ws.Range("C25").WrapText = $true
$ws.Range("C25").VerticalAlignment = -4108

$ws.Range("D25").Value = "2h"
$ws.Range("D25").WrapText = $true
$ws.Range("D25").VerticalAlignment = -4108
$ws.Range("D25").HorizontalAlignment = -4108

$ws.Range("E25").Value = "Team manager "
$ws.Range("E25").WrapText = $true
$ws.Range("E25").VerticalAlignment = -4108

# --- Row 26 ---
$ws.Range("B26").Value = "truck's driver"
$ws.Range("B26").WrapText = $true
$ws.Range("B26").VerticalAlignment = -4108

$ws.Range("C26").Value = "Get notification of changing route "
$ws.Range("C26").WrapText = $true
$ws.Range("C26").VerticalAlignment = -4108

$ws.Range("D26").Value = "1h"
$ws.Range("D26").WrapText = $true
$ws.Range("D26").VerticalAlignment = -4108
$ws.Range("D26").HorizontalAlignment = -4108

$ws.Range("E26").Value = "All truck's driver"
$ws.Range("E26").WrapText = $true
$ws.Range("E26").VerticalAlignment = -4108

# --- Selection cosmetics ---
$ws.Range("F26").Select()
